$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.48757729316263
$ws.Range("D2").Value = 0.006833957729854975
$ws.Range("E2").Value = 1.155934970043793
$ws.Range("F2").Value = 1.336185654970144

$ws.Range("B3").Value = 11.48605387034872
$ws.Range("D3").Value = 0.006823036354184879
$ws.Range("E3").Value = 1.15408766566219
$ws.Range("F3").Value = 1.331918340033602

$ws.Range("B4").Value = 11.36876089678898
$ws.Range("D4").Value = 0.006847703136742046
$ws.Range("E4").Value = 1.15825994733021
$ws.Range("F4").Value = 1.34156610558938

$ws.Range("B5").Value = 11.29719582210498
$ws.Range("D5").Value = 0.007106744882627307
$ws.Range("E5").Value = 1.202075760158821
$ws.Range("F5").Value = 1.444986133161408

$ws.Range("B6").Value = 11.36718673505485
$ws.Range("D6").Value = 0.006839232491263472
$ws.Range("E6").Value = 1.156827173567979
$ws.Range("F6").Value = 1.33824910950528
